# "Source:" / "La Contraloría General de la República" / <URL> used to be three
# consecutive rows (36-38) with the URL row (38) carrying a live hyperlink, and
# a trailing blank "spacer" row (39) after them.
#
# The edit turns that into: a blank spacer row right after "Source:" (so it
# reads "Source:" / blank / "La Contraloría..." / blank / URL), and the URL
# text is demoted from a hyperlink to plain italic ("source"-styled) text -
# the <hyperlinks> part of the sheet goes away entirely.
#
# We get there by inserting a fresh blank row at 37 (pushing the old rows
# 37-39 down to 38-40, each keeping its original formatting), dropping the
# hyperlink object, and then fixing up the one row (the old hyperlink row,
# now 39) whose style needs to change from "HyperLink" back to "source" - the
# cleanest way we found to do that reliably is to delete that single row and
# re-insert a blank one, which picks up the "source" style from row 38 above
# it, exactly like the very first insert did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 37; rows 37-39 ("La Contraloría...", the URL with
# its hyperlink, and the trailing blank row) all shift down to 38-40 and keep
# their existing formatting.
$ws.Rows("37:37").Insert()

# The URL text now lives at A39 (still carrying the old hyperlink + its
# "HyperLink" styling) - grab its text before anything else changes.
$url = $ws.Range("A39").Value2

# Remove the hyperlink itself.
$ws.Hyperlinks.Delete()

# A39 still has the old "HyperLink" look; deleting the row and inserting a
# fresh blank one in its place makes it inherit A38's "source" style instead,
# the same way the blank row at 37 already did.
$ws.Rows("39:39").Delete()
$ws.Rows("39:39").Insert()

# Finally, put the URL text - now plain, "source"-styled text - onto the new
# last row, 40.
$ws.Range("A40").Value2 = $url
